# Reproduce the target edit: a new (blank) column is inserted before column A
# and a new (blank) row is inserted above row 1, shifting the existing
# A1:C3 table down/right to B2:D4. The relocated table is then re-populated
# with new host-list data, and row 3 / row 4 only use columns B and C
# (column D is left empty on those rows). Finally the active selection is
# moved to C5, matching the saved workbook's cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before A and a blank row above row 1. This shifts
# the old A1:C3 content (and the column-A style/width) to B2:D4.
[void]$ws.Columns("A").Insert()
[void]$ws.Rows("1").Insert()

# Row 2 (was row 1): 123.57.56.121 / abc / alias
$ws.Range("B2").Value = "123.57.56.121"
$ws.Range("C2").Value = "abc"
$ws.Range("D2").Value = "alias"

# Row 3 (was row 2): 123.57.56.1 / liujiashuai   (column D left blank)
$ws.Range("B3").Value = "123.57.56.1"
$ws.Range("C3").Value = "liujiashuai"
$ws.Range("D3").ClearContents()

# Row 4 (was row 3): 123.57.56.121 / MARS        (column D left blank)
$ws.Range("B4").Value = "123.57.56.121"
$ws.Range("C4").Value = "MARS"
$ws.Range("D4").ClearContents()

# Match the saved cursor/selection position.
[void]$ws.Range("C5").Select()
